$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the cadence from monthly back to weekly starting at row 61 (B61).
$ws.Cells.Item(61, 2).Value = 44381
$ws.Cells.Item(62, 2).Value = 44388

# Copy formatting (borders/number formats/font) from the last existing
# data row down onto the newly appended rows before filling them in.
$src = $ws.Range("A61:B61")
$dst = $ws.Range("A63:B70")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats

# Append the new weekly rows (index 61-68, dates 7 days apart).
$newValues = @(
    @(61, 44395),
    @(62, 44402),
    @(63, 44409),
    @(64, 44416),
    @(65, 44423),
    @(66, 44430),
    @(67, 44437),
    @(68, 44444)
)

$row = 63
foreach ($pair in $newValues) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
